# Update the "New Rules" sheet:
#  1. Fix a typo: C34 was "MINOE" -> should be "MINOR".
#  2. Append five new rule rows (35-39) describing new PMD-style rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Rules")

# --- 1. Typo fix -----------------------------------------------------
$ws.Range("C34").Value = "MINOR"

# --- 2. New rule rows --------------------------------------------------
# Row 35 - queueNamingConvention
$ws.Range("A35").Value = "queueNamingConvention"
$ws.Range("B35").Value = "MQ Nodes"
$ws.Range("C35").Value = "MINOR"
$ws.Range("D35").Value = "15 min"
$ws.Range("E35").Value = "Standards"
$ws.Range("F35").Value = "MQ nodes should access alias queues. The naming convetion for alias queues is '^[A-Za-z0-9_]+\.[A-Za-z0-9_]+\.[A-Za-z0-9_]+\.(AI|AO)$'."
$ws.Rows.Item(35).RowHeight = 30

# Row 36 - InconsistentRouteNode
$ws.Range("A36").Value = "InconsistentRouteNode"
$ws.Range("B36").Value = "Route Node"
$ws.Range("C36").Value = "CRITICAL"
$ws.Range("D36").Value = "1 h"
$ws.Range("E36").Value = "Correctness"
$ws.Range("F36").Value = "All the terminals specified in the filter table should be connected, otherwise it may cause abnormal termination and message might be lost."

# Row 37 - DeprecatedNodeCheck
$ws.Range("A37").Value = "DeprecatedNodeCheck"
$ws.Range("B37").Value = "Check Node"
$ws.Range("C37").Value = "MAJOR"
$ws.Range("D37").Value = "5 min"
$ws.Range("E37").Value = "Deprecated"
$ws.Range("F37").Value = "Usage of deprecated nodes in the message flow is discouraged."

# Row 38 - DSNWithoutDBCall
$ws.Range("A38").Value = "DSNWithoutDBCall"
$ws.Range("B38").Value = "Compute Node"
$ws.Range("C38").Value = "MAJOR"
$ws.Range("D38").Value = "10 min"
$ws.Range("E38").Value = "Correctness,`nPerformance"
$ws.Range("F38").Value = "Data Source should not be specified if compute node is not interacting with Database."
$ws.Range("E38").WrapText = $true
$ws.Rows.Item(38).RowHeight = 30

# Row 39 - MavenProjectNamingConventions
$ws.Range("A39").Value = "MavenProjectNamingConventions"
$ws.Range("B39").Value = "Project"
$ws.Range("C39").Value = "MINOR"
$ws.Range("D39").Value = "30 min"
$ws.Range("E39").Value = "Standards"
$ws.Range("F39").Value = "The maven project artifacts and modules should follow Naming conventions."

# --- Formatting for column A/B/F to match existing style (wrap text on F) ---
$ws.Range("F35:F39").WrapText = $true
$ws.Range("A35:A39").Font.Bold = $false

# --- View: move selection to C40 and drop the old frozen top-left cell ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C40").Select()
